$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.126.03"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.051.94"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'248.36"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").Value = "'57.38"
$ws.Range("E7").Value = "  -3.66%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").Value = "'0.0775"
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'15.99"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "'0.878"
$ws.Range("E13").Value = "  +7.18%  "
$ws.Range("D14").Value = "2.348.43"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "2.051.77"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "'17.90"
$ws.Range("E17").Value = "  +14.12%  "
$ws.Range("D18").Value = "37.162.07"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'74.83"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("D21").Value = "'5.37"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").Value = "'237.25"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").Value = "'9.49"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = "  -5.23%  "
$ws.Range("D27").Value = "'169.23"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'20.04"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "'4.81"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "'0.0618"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'0.0896"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("D37").Value = "'1.78"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.34"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Value = "'3.21"
$ws.Range("E39").Value = "  +13.54%  "
$ws.Range("D40").Value = "'5.17"
$ws.Range("E40").Value = "  +14.85%  "
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").Value = "'17.28"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "'0.0938"
$ws.Range("E44").Value = "  -21.55%  "
$ws.Range("D45").Value = "'95.87"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "1.273.47"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").Value = "'6.82"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").Value = "2.232.32"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "'43.57"
$ws.Range("E51").Value = "  -1.26%  "
